$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new data row 33 following the same pattern as the preceding rows
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 110032
$ws.Range("C33").Value = 10032
$ws.Range("D33").Value = "eng"
$ws.Range("E33").Value = $true
$ws.Range("F33").Value = "superadmin"
$ws.Range("G33").Value = "now()"

# Update selection to match the new active cell shown in the saved workbook
$ws.Range("C30").Select()
